# This workbook tracks weekly Plátano (banana) price observations for the
# "Terminal Hortofrutícola Agro Chillán" market. Each reporting week adds a
# pair of rows (one for "Pintón" quality, one for "Primera Pintón" quality)
# at the top of the historical data block (row 236), pushing the rest of the
# history down by two rows. This script reproduces that weekly append.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh rows at the top of the data block (row 236/237). Excel
# shifts every existing row from 236 downward to 238.. automatically,
# carrying along all of their values/formatting - exactly matching the
# "old row N becomes new row N+2" relationship seen between the two
# workbook revisions.
$ws.Rows("236:237").Insert()

# --- New row 236: "Pintón" quality, week of 2021-09-27 -----------------
$ws.Range("A236").Value2 = 7
$ws.Range("B236").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C236").Value = "Ñuble"
$ws.Range("D236").Value2 = 44466
$ws.Range("E236").Value2 = 16
$ws.Range("F236").Value = "Fruta"
$ws.Range("G236").Value2 = 100108
$ws.Range("H236").Value = "Tropicales y subtropicales"
$ws.Range("I236").Value2 = 100108006
$ws.Range("J236").Value = "Plátano"
$ws.Range("K236").Value = "Sin especificar"
$ws.Range("L236").Value = "Pintón"
$ws.Range("M236").Value2 = 150
$ws.Range("N236").Value2 = 14000
$ws.Range("O236").Value2 = 14000
$ws.Range("P236").Value2 = 14000
$ws.Range("Q236").Value = "$/caja 20 kilos"
$ws.Range("R236").Value = "Ecuador"
$ws.Range("S236").Value2 = 700
$ws.Range("T236").Value2 = 20

# --- New row 237: "Primera Pintón" quality, same week ------------------
$ws.Range("A237").Value2 = 7
$ws.Range("B237").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C237").Value = "Ñuble"
$ws.Range("D237").Value2 = 44466
$ws.Range("E237").Value2 = 16
$ws.Range("F237").Value = "Fruta"
$ws.Range("G237").Value2 = 100108
$ws.Range("H237").Value = "Tropicales y subtropicales"
$ws.Range("I237").Value2 = 100108006
$ws.Range("J237").Value = "Plátano"
$ws.Range("K237").Value = "Sin especificar"
$ws.Range("L237").Value = "Primera Pintón"
$ws.Range("M237").Value2 = 400
$ws.Range("N237").Value2 = 15000
$ws.Range("O237").Value2 = 16000
$ws.Range("P237").Value2 = 15500
$ws.Range("Q237").Value = "$/caja 20 kilos"
$ws.Range("R237").Value = "Ecuador"
$ws.Range("S237").Value2 = 775
$ws.Range("T237").Value2 = 20

# Keep the D column's date number format/style consistent with every other
# date cell in the column (style index 2 in the original workbook).
$ws.Range("D236:D237").Style = $ws.Range("D238").Style
